{"js": "// Change location of README\n// 1) \"...customer management app in use, or are looking...\" -> remove the comma\n//    after \"in use\" (\"in use, or\" -> \"in use or\").\n// 2) \"...delivery status of orders. Additionally it lets the use...\" -> add a\n//    comma after \"Additionally\" (\"Additionally it lets\" -> \"Additionally, it lets\").\n\nconst body = context.document.body;\n\nconst commaResults = body.search(\"in use,\", { matchCase: true, matchWholeWord: false });\ncommaResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < commaResults.items.length; i++) {\n  commaResults.items[i].insertText(\"in use\", \"Replace\");\n}\nawait context.sync();\n\nconst additionallyResults = body.search(\"Additionally it lets\", { matchCase: true, matchWholeWord: false });\nadditionallyResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < additionallyResults.items.length; i++) {\n  additionallyResults.items[i].insertText(\"Additionally, it lets\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Change location of README\n# 1) \"...customer management app in use, or are looking...\" -> remove the comma\n#    after \"in use\" (\"in use, or\" -> \"in use or\").\n# 2) \"...delivery status of orders. Additionally it lets the use...\" -> add a\n#    comma after \"Additionally\" (\"Additionally it lets\" -> \"Additionally, it lets\").\n\n$d = $word.ActiveDocument\n\n$range1 = $d.Content\n$range1.Find.Execute(\"in use,\", $false, $false, $false, $false, $false, $true, 1, $false, \"in use\", 2)\n\n$range2 = $d.Content\n$range2.Find.Execute(\"Additionally it lets\", $false, $false, $false, $false, $false, $true, 1, $false, \"Additionally, it lets\", 2)\n"}
